# RAD Extension Payment Type test-data update
# - Clears the "Execute" (Y) flag for every scenario row except the
#   "Extension Payments" row (row 4), which is the scenario now being run.
# - Stamps the "Extension Payments" row's Date column with the new run time.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Clear()
$ws.Range("C3").Clear()
$ws.Range("C5").Clear()
$ws.Range("C6").Clear()
$ws.Range("C7").Clear()

$ws.Range("B4").Value = "Wed Mar 20 23:05:16 EDT 2024"
